$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 32 (shifts rows 32..332 down to 33..333, carrying
# formatting/heights with them, matching the table's existing row styles).
$ws.Rows("32:32").Insert()

# Populate the new row 32 with the new "strWindowPos" string resource.
$ws.Range("B32").Value = "localization\strings"
$ws.Range("C32").Value = "strWindowPos"
$ws.Range("D32").Value = "In ""settings"" form, tab ""User interface"""
$ws.Range("E32").Value = "Remember window position and size on startup"

# The existing "strChkDlgPath" row (now row 25) gains a Comment value.
$ws.Range("D25").Value = "In ""settings"" form, tab ""User interface"""

# Extend the "Tabla13" table so it covers the newly inserted row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:F204"))

# Column D was widened (to fit the new/longer Comment text).
$ws.Columns("D").ColumnWidth = 34.8
